$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodeSchemes")

# Remember column A's width so the freshly inserted column can match it
# (mirrors Excel's normal "insert column" behaviour of carrying over the
# neighbouring column's width).
$origWidth = $ws.Columns.Item(1).ColumnWidth

# Insert a new column before column B (ID), shifting existing columns right.
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(2).ColumnWidth = $origWidth

# Populate the new ORGANIZATION column.
$ws.Cells.Item(1,2).Value = "ORGANIZATION"
$ws.Cells.Item(2,2).Value = "74a41211-8c99-4835-a519-7a61612b1098"

# Update the defined name "yti" so its range keeps covering the header/data
# rows now that an extra column has been added (X -> Y).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "CodeSchemes!yti") {
        $n.RefersTo = "=CodeSchemes!`$A`$1:`$Y`$2"
    }
}
